# Generate Report for Handoff
# Updates the Priority column for the "Ready for handoff" rows to "ht"
# on both the zh-cn and de-de worksheets, and refreshes the associated
# handoff timestamps (Overview "Latest HO Xliff Generate Date" and the
# per-language "Latest Handoff Datetime") to reflect the regenerated
# handoff xliff files.

$wb = $excel.ActiveWorkbook

$rows = @(7, 10, 11, 12, 13, 14)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsOverview = $wb.Worksheets.Item("Overview")

foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"

    $wsZhCn.Range("H$r").Value = "2016-08-17 08:21:19"
    $wsDeDe.Range("H$r").Value = "2016-08-17 08:21:24"
    $wsOverview.Range("G$r").Value = "2016-08-17 08:21:24"
}
